$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brix_Gel_Stab")

# Shift rows 4..10 down by one (process bottom-up so we don't clobber data
# before it's copied) to make room for the new "Schräge" row at row 4.
for ($r = 10; $r -ge 4; $r--) {
    $dest = $r + 1
    for ($c = 1; $c -le 5; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $destCell = $ws.Cells.Item($dest, $c)
        $destCell.Value = $srcCell.Value2
    }
}

# New row 4: "Schräge" option (roof/stair slope factor)
$ws.Cells.Item(4, 1).Value = "Auswahl"
$ws.Cells.Item(4, 2).Value = "Schräge"
$ws.Cells.Item(4, 3).Value = "F_Schräg"
$ws.Cells.Item(4, 4).Value = "´---:0, bis 6°:32, bis 35°:60"

# Fix typo in the "Montage Steher" options string (was row 6, now row 7):
# leading backtick (`) corrected to acute accent (´)
$ws.Cells.Item(7, 4).Value = "´---:0, Aufsatz:125, Seite:161"

# Update final price formula (was row 10, now row 11) to include the new
# F_Schräg term
$ws.Cells.Item(11, 5).Value = "((P_Modell * L * F_Faktor * P_Handlauf) + ((math.ceil(L/1.3)+1) * P_Steher * F_Faktor) + (Ecken * 95) + (L * P_Arbeit) + (L * F_Schräg)) * ( 1 - (p_rabatt / 100))"

$ws.Range("E11").Select() | Out-Null
